$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $val) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextCell 'D2' '27.461.34'
Set-TextCell 'E2' '  +2.07%  '
Set-TextCell 'D3' '1.838.51'
Set-TextCell 'E3' '  +1.35%  '
Set-TextCell 'D4' '1.012'
Set-TextCell 'E4' '  +1.02%  '
Set-TextCell 'D5' '314.39'
Set-TextCell 'E5' '  +1.70%  '
Set-TextCell 'D6' '1.011'
Set-TextCell 'E6' '  +0.96%  '
Set-TextCell 'D7' '0.4742'
Set-TextCell 'E7' '  +1.70%  '
Set-TextCell 'D8' '0.3693'
Set-TextCell 'E8' '  +0.89%  '
Set-TextCell 'D9' '0.07460'
Set-TextCell 'E9' '  +1.47%  '
Set-TextCell 'D10' '0.8850'
Set-TextCell 'E10' '  +1.88%  '
Set-TextCell 'D11' '20.45'
Set-TextCell 'E11' '  +0.70%  '
Set-TextCell 'D12' '1.926.64'
Set-TextCell 'E12' '  +6.65%  '
Set-TextCell 'D13' '0.07332'
Set-TextCell 'E13' '  +3.54%  '
Set-TextCell 'D14' '5.452'
Set-TextCell 'E14' '  +1.20%  '
Set-TextCell 'D15' '93.30'
Set-TextCell 'E15' '  +1.83%  '
Set-TextCell 'D16' '6.584'
Set-TextCell 'E16' '  +1.04%  '
Set-TextCell 'D17' '1.010'
Set-TextCell 'E17' '  +0.71%  '
Set-TextCell 'D18' '0.000008817'
Set-TextCell 'E18' '  +1.27%  '
Set-TextCell 'E19' '  +0.99%  '
Set-TextCell 'D20' '14.81'
Set-TextCell 'E20' '  +1.07%  '
Set-TextCell 'D21' '27.505.82'
Set-TextCell 'D22' '5.327'
Set-TextCell 'E22' '  +0.60%  '
Set-TextCell 'D23' '10.70'
Set-TextCell 'E23' '  +0.66%  '
Set-TextCell 'D24' '2.142.61'
Set-TextCell 'E24' '  +5.20%  '
Set-TextCell 'D25' '1.908'
Set-TextCell 'E25' '  +0.69%  '
Set-TextCell 'D26' '152.19'
Set-TextCell 'E26' '  +0.89%  '
Set-TextCell 'D27' '18.64'
Set-TextCell 'E27' '  +1.50%  '
Set-TextCell 'D28' '2.149'
Set-TextCell 'E28' '  -0.13%  '
Set-TextCell 'D29' '5.254'
Set-TextCell 'E29' '  -0.22%  '
Set-TextCell 'D30' '117.99'
Set-TextCell 'E30' '  +2.26%  '
Set-TextCell 'D31' '0.08973'
Set-TextCell 'E31' '  +0.84%  '
Set-TextCell 'D32' '0.7584'
Set-TextCell 'E32' '  +0.52%  '
Set-TextCell 'D33' '1.184'
Set-TextCell 'E33' '  +2.38%  '
Set-TextCell 'D34' '4.563'
Set-TextCell 'E34' '  +1.66%  '
Set-TextCell 'D35' '2.951'
Set-TextCell 'E35' '  +1.35%  '
Set-TextCell 'D36' '1.013'
Set-TextCell 'E36' '  +1.14%  '
Set-TextCell 'D37' '1.105'
Set-TextCell 'E37' '  +1.96%  '
Set-TextCell 'D38' '0.05331'
Set-TextCell 'E38' '  +1.06%  '
Set-TextCell 'D39' '0.01957'
Set-TextCell 'E39' '  +0.47%  '
Set-TextCell 'D40' '2.998'
Set-TextCell 'E40' '  +1.37%  '
Set-TextCell 'D41' '7.347'
Set-TextCell 'E41' '  +1.71%  '
Set-TextCell 'D42' '2.411'
Set-TextCell 'E42' '  +5.68%  '
Set-TextCell 'E43' '  +0.71%  '
Set-TextCell 'D44' '0.1660'
Set-TextCell 'E44' '  +0.58%  '
Set-TextCell 'D45' '8.524'
Set-TextCell 'E45' '  +1.28%  '
Set-TextCell 'D46' '0.4922'
Set-TextCell 'E46' '  +1.03%  '
Set-TextCell 'D47' '10.53'
Set-TextCell 'E47' '  +1.67%  '
Set-TextCell 'B48' 'Quant'
Set-TextCell 'C48' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextCell 'D48' '105.22'
Set-TextCell 'E48' '  +2.01%  '
Set-TextCell 'B49' 'PaxDollar'
Set-TextCell 'C49' 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextCell 'D49' '1.012'
Set-TextCell 'E49' '  +1.10%  '
Set-TextCell 'E50' '  +1.22%  '
Set-TextCell 'D51' '0.06317'
